$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 186-223: shift each row's date/volume/price data down from the previous row,
# and seed a brand-new data point into row 186.
$ws.Cells.Item(186, 4).Value = 44511
$ws.Cells.Item(186, 10).Value = 40
$ws.Cells.Item(186, 11).Value = 4000
$ws.Cells.Item(186, 12).Value = 4000
$ws.Cells.Item(186, 13).Value = 4000
$ws.Cells.Item(186, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(186, 16).Value = 1333

$ws.Cells.Item(187, 4).Value = 44306
$ws.Cells.Item(187, 10).Value = 35
$ws.Cells.Item(187, 11).Value = 3000
$ws.Cells.Item(187, 12).Value = 3000
$ws.Cells.Item(187, 13).Value = 3000
$ws.Cells.Item(187, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(187, 16).Value = 1000

$ws.Cells.Item(188, 4).Value = 44211
$ws.Cells.Item(188, 10).Value = 65
$ws.Cells.Item(188, 11).Value = 3000
$ws.Cells.Item(188, 12).Value = 3000
$ws.Cells.Item(188, 13).Value = 3000
$ws.Cells.Item(188, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(188, 16).Value = 1000

$ws.Cells.Item(189, 4).Value = 44469
$ws.Cells.Item(189, 10).Value = 60
$ws.Cells.Item(189, 11).Value = 4000
$ws.Cells.Item(189, 12).Value = 5000
$ws.Cells.Item(189, 13).Value = 4500
$ws.Cells.Item(189, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(189, 16).Value = 1500

$ws.Cells.Item(190, 4).Value = 44215
$ws.Cells.Item(190, 10).Value = 40
$ws.Cells.Item(190, 11).Value = 3000
$ws.Cells.Item(190, 12).Value = 3000
$ws.Cells.Item(190, 13).Value = 3000
$ws.Cells.Item(190, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(190, 16).Value = 1000

$ws.Cells.Item(191, 4).Value = 44186
$ws.Cells.Item(191, 10).Value = 40
$ws.Cells.Item(191, 11).Value = 4000
$ws.Cells.Item(191, 12).Value = 4000
$ws.Cells.Item(191, 13).Value = 4000
$ws.Cells.Item(191, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(191, 16).Value = 1333

$ws.Cells.Item(192, 4).Value = 44504
$ws.Cells.Item(192, 10).Value = 45
$ws.Cells.Item(192, 11).Value = 5000
$ws.Cells.Item(192, 12).Value = 5000
$ws.Cells.Item(192, 13).Value = 5000
$ws.Cells.Item(192, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(192, 16).Value = 1667

$ws.Cells.Item(193, 4).Value = 44246
$ws.Cells.Item(193, 10).Value = 65
$ws.Cells.Item(193, 11).Value = 4000
$ws.Cells.Item(193, 12).Value = 4000
$ws.Cells.Item(193, 13).Value = 4000
$ws.Cells.Item(193, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(193, 16).Value = 1333

$ws.Cells.Item(194, 4).Value = 44505
$ws.Cells.Item(194, 10).Value = 45
$ws.Cells.Item(194, 11).Value = 5000
$ws.Cells.Item(194, 12).Value = 5000
$ws.Cells.Item(194, 13).Value = 5000
$ws.Cells.Item(194, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(194, 16).Value = 1667

$ws.Cells.Item(195, 4).Value = 44487
$ws.Cells.Item(195, 10).Value = 105
$ws.Cells.Item(195, 11).Value = 4000
$ws.Cells.Item(195, 12).Value = 4500
$ws.Cells.Item(195, 13).Value = 4310
$ws.Cells.Item(195, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(195, 16).Value = 1437

$ws.Cells.Item(196, 4).Value = 44487
$ws.Cells.Item(196, 10).Value = 65
$ws.Cells.Item(196, 11).Value = 2500
$ws.Cells.Item(196, 12).Value = 2500
$ws.Cells.Item(196, 13).Value = 2500
$ws.Cells.Item(196, 15).Value = "Región Metropolitana"
$ws.Cells.Item(196, 16).Value = 833

$ws.Cells.Item(197, 4).Value = 44425
$ws.Cells.Item(197, 10).Value = 40
$ws.Cells.Item(197, 11).Value = 4000
$ws.Cells.Item(197, 12).Value = 4000
$ws.Cells.Item(197, 13).Value = 4000
$ws.Cells.Item(197, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(197, 16).Value = 1333

$ws.Cells.Item(198, 4).Value = 44425
$ws.Cells.Item(198, 10).Value = 10
$ws.Cells.Item(198, 11).Value = 3300
$ws.Cells.Item(198, 12).Value = 3300
$ws.Cells.Item(198, 13).Value = 3300
$ws.Cells.Item(198, 15).Value = "Región Metropolitana"
$ws.Cells.Item(198, 16).Value = 1100

$ws.Cells.Item(199, 4).Value = 44343
$ws.Cells.Item(199, 10).Value = 30
$ws.Cells.Item(199, 11).Value = 4000
$ws.Cells.Item(199, 12).Value = 4000
$ws.Cells.Item(199, 13).Value = 4000
$ws.Cells.Item(199, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(199, 16).Value = 1333

$ws.Cells.Item(200, 4).Value = 44370
$ws.Cells.Item(200, 10).Value = 20
$ws.Cells.Item(200, 11).Value = 4000
$ws.Cells.Item(200, 12).Value = 4000
$ws.Cells.Item(200, 13).Value = 4000
$ws.Cells.Item(200, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(200, 16).Value = 1333

$ws.Cells.Item(201, 4).Value = 44449
$ws.Cells.Item(201, 10).Value = 65
$ws.Cells.Item(201, 11).Value = 4000
$ws.Cells.Item(201, 12).Value = 4000
$ws.Cells.Item(201, 13).Value = 4000
$ws.Cells.Item(201, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(201, 16).Value = 1333

$ws.Cells.Item(202, 4).Value = 44168
$ws.Cells.Item(202, 10).Value = 125
$ws.Cells.Item(202, 11).Value = 4500
$ws.Cells.Item(202, 12).Value = 5000
$ws.Cells.Item(202, 13).Value = 4740
$ws.Cells.Item(202, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(202, 16).Value = 1580

$ws.Cells.Item(203, 4).Value = 44175
$ws.Cells.Item(203, 10).Value = 50
$ws.Cells.Item(203, 11).Value = 4000
$ws.Cells.Item(203, 12).Value = 4000
$ws.Cells.Item(203, 13).Value = 4000
$ws.Cells.Item(203, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(203, 16).Value = 1333

$ws.Cells.Item(204, 4).Value = 44392
$ws.Cells.Item(204, 10).Value = 35
$ws.Cells.Item(204, 11).Value = 4000
$ws.Cells.Item(204, 12).Value = 4000
$ws.Cells.Item(204, 13).Value = 4000
$ws.Cells.Item(204, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(204, 16).Value = 1333

$ws.Cells.Item(205, 4).Value = 44286
$ws.Cells.Item(205, 10).Value = 20
$ws.Cells.Item(205, 11).Value = 3000
$ws.Cells.Item(205, 12).Value = 3000
$ws.Cells.Item(205, 13).Value = 3000
$ws.Cells.Item(205, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(205, 16).Value = 1000

$ws.Cells.Item(206, 4).Value = 44473
$ws.Cells.Item(206, 10).Value = 50
$ws.Cells.Item(206, 11).Value = 4000
$ws.Cells.Item(206, 12).Value = 4000
$ws.Cells.Item(206, 13).Value = 4000
$ws.Cells.Item(206, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(206, 16).Value = 1333

$ws.Cells.Item(207, 4).Value = 44400
$ws.Cells.Item(207, 10).Value = 40
$ws.Cells.Item(207, 11).Value = 4000
$ws.Cells.Item(207, 12).Value = 5000
$ws.Cells.Item(207, 13).Value = 4500
$ws.Cells.Item(207, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(207, 16).Value = 1500

$ws.Cells.Item(208, 4).Value = 44484
$ws.Cells.Item(208, 10).Value = 40
$ws.Cells.Item(208, 11).Value = 4000
$ws.Cells.Item(208, 12).Value = 4500
$ws.Cells.Item(208, 13).Value = 4250
$ws.Cells.Item(208, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(208, 16).Value = 1417

$ws.Cells.Item(209, 4).Value = 44181
$ws.Cells.Item(209, 10).Value = 55
$ws.Cells.Item(209, 11).Value = 5000
$ws.Cells.Item(209, 12).Value = 5000
$ws.Cells.Item(209, 13).Value = 5000
$ws.Cells.Item(209, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(209, 16).Value = 1667

$ws.Cells.Item(210, 4).Value = 44342
$ws.Cells.Item(210, 10).Value = 20
$ws.Cells.Item(210, 11).Value = 4000
$ws.Cells.Item(210, 12).Value = 4000
$ws.Cells.Item(210, 13).Value = 4000
$ws.Cells.Item(210, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(210, 16).Value = 1333

$ws.Cells.Item(211, 4).Value = 44328
$ws.Cells.Item(211, 10).Value = 55
$ws.Cells.Item(211, 11).Value = 4000
$ws.Cells.Item(211, 12).Value = 4000
$ws.Cells.Item(211, 13).Value = 4000
$ws.Cells.Item(211, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(211, 16).Value = 1333

$ws.Cells.Item(212, 4).Value = 44301
$ws.Cells.Item(212, 10).Value = 50
$ws.Cells.Item(212, 11).Value = 3000
$ws.Cells.Item(212, 12).Value = 3000
$ws.Cells.Item(212, 13).Value = 3000
$ws.Cells.Item(212, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(212, 16).Value = 1000

$ws.Cells.Item(213, 4).Value = 44330
$ws.Cells.Item(213, 10).Value = 20
$ws.Cells.Item(213, 11).Value = 4000
$ws.Cells.Item(213, 12).Value = 4000
$ws.Cells.Item(213, 13).Value = 4000
$ws.Cells.Item(213, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(213, 16).Value = 1333

$ws.Cells.Item(214, 4).Value = 44270
$ws.Cells.Item(214, 10).Value = 30
$ws.Cells.Item(214, 11).Value = 3000
$ws.Cells.Item(214, 12).Value = 3000
$ws.Cells.Item(214, 13).Value = 3000
$ws.Cells.Item(214, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(214, 16).Value = 1000

$ws.Cells.Item(215, 4).Value = 44217
$ws.Cells.Item(215, 10).Value = 50
$ws.Cells.Item(215, 11).Value = 3000
$ws.Cells.Item(215, 12).Value = 3000
$ws.Cells.Item(215, 13).Value = 3000
$ws.Cells.Item(215, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(215, 16).Value = 1000

$ws.Cells.Item(216, 4).Value = 44509
$ws.Cells.Item(216, 10).Value = 20
$ws.Cells.Item(216, 11).Value = 4000
$ws.Cells.Item(216, 12).Value = 4000
$ws.Cells.Item(216, 13).Value = 4000
$ws.Cells.Item(216, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(216, 16).Value = 1333

$ws.Cells.Item(217, 4).Value = 44383
$ws.Cells.Item(217, 10).Value = 30
$ws.Cells.Item(217, 11).Value = 4000
$ws.Cells.Item(217, 12).Value = 4000
$ws.Cells.Item(217, 13).Value = 4000
$ws.Cells.Item(217, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(217, 16).Value = 1333

$ws.Cells.Item(218, 4).Value = 44307
$ws.Cells.Item(218, 10).Value = 35
$ws.Cells.Item(218, 11).Value = 3000
$ws.Cells.Item(218, 12).Value = 3000
$ws.Cells.Item(218, 13).Value = 3000
$ws.Cells.Item(218, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(218, 16).Value = 1000

$ws.Cells.Item(219, 4).Value = 44273
$ws.Cells.Item(219, 10).Value = 50
$ws.Cells.Item(219, 11).Value = 3000
$ws.Cells.Item(219, 12).Value = 3000
$ws.Cells.Item(219, 13).Value = 3000
$ws.Cells.Item(219, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(219, 16).Value = 1000

$ws.Cells.Item(220, 4).Value = 44433
$ws.Cells.Item(220, 10).Value = 55
$ws.Cells.Item(220, 11).Value = 4000
$ws.Cells.Item(220, 12).Value = 4000
$ws.Cells.Item(220, 13).Value = 4000
$ws.Cells.Item(220, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(220, 16).Value = 1333

$ws.Cells.Item(221, 4).Value = 44302
$ws.Cells.Item(221, 10).Value = 40
$ws.Cells.Item(221, 11).Value = 3000
$ws.Cells.Item(221, 12).Value = 3000
$ws.Cells.Item(221, 13).Value = 3000
$ws.Cells.Item(221, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(221, 16).Value = 1000

$ws.Cells.Item(222, 4).Value = 44179
$ws.Cells.Item(222, 10).Value = 20
$ws.Cells.Item(222, 11).Value = 5000
$ws.Cells.Item(222, 12).Value = 5000
$ws.Cells.Item(222, 13).Value = 5000
$ws.Cells.Item(222, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(222, 16).Value = 1667

$ws.Cells.Item(223, 4).Value = 44491
$ws.Cells.Item(223, 10).Value = 65
$ws.Cells.Item(223, 11).Value = 4000
$ws.Cells.Item(223, 12).Value = 4000
$ws.Cells.Item(223, 13).Value = 4000
$ws.Cells.Item(223, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(223, 16).Value = 1333

# Row 224 is new: append the data point that previously lived at the end of the series.
$ws.Cells.Item(224, 1).Value = 10
$ws.Cells.Item(224, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(224, 3).Value = "La Araucanía"
$ws.Cells.Item(224, 4).Value = 44491
$ws.Cells.Item(224, 5).Value = 9
$ws.Cells.Item(224, 6).Value = 100112044
$ws.Cells.Item(224, 7).Value = "Perejil"
$ws.Cells.Item(224, 8).Value = "Sin especificar"
$ws.Cells.Item(224, 9).Value = "Primera"
$ws.Cells.Item(224, 10).Value = 40
$ws.Cells.Item(224, 11).Value = 5000
$ws.Cells.Item(224, 12).Value = 5000
$ws.Cells.Item(224, 13).Value = 5000
$ws.Cells.Item(224, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(224, 15).Value = "Región del Maule"
$ws.Cells.Item(224, 16).Value = 1667
$ws.Cells.Item(224, 17).Value = 3
$ws.Cells.Item(224, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(224, 4).NumberFormat = $ws.Cells.Item(223, 4).NumberFormat
